$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 6 ("Extension.valueCodeableConcept" -> "Extension.value[x]") edits
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Extension.value[x]"
$ws.Range("B6").ClearContents()
$ws.Range("E6").Value = "0"
$ws.Range("K6").Value = "Value of extension"
$ws.Range("W6").ClearContents()
$ws.Range("X6").ClearContents()
$ws.Range("Y6").ClearContents()
$ws.Range("AA6").Value = 'type:$this}' + "`n"
$ws.Range("AB6").ClearContents()
$ws.Range("AD6").Value = "closed"

# ---------------------------------------------------------------------------
# 2. New row 7 - the "valueCodeableConcept" slice of Extension.value[x]
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Extension.value[x]"
$ws.Range("B7").Value = "valueCodeableConcept"
$ws.Range("E7").Value = "1"
$ws.Range("F7").Value = "1"
$ws.Range("J7").Value = "CodeableConcept" + "`n"
$ws.Range("K7").Value = "This records the outcome of an Out-Patient Attendance Consultant."
$ws.Range("L7").Value = 'Value of extension - may be a resource or one of a constrained set of the data types (see Extensibility in the spec for list).'
$ws.Range("W7").Value = "required"
$ws.Range("X7").Value = "This records the outcome of an Out-Patient Attendance Consultant."
$ws.Range("Y7").Value = "https://fhir.hl7.org.uk/STU3/ValueSet/CareConnect-OutcomeOfAttendance-1"
$ws.Range("AE7").Value = "Extension.value[x]"
$ws.Range("AF7").Value = "0"
$ws.Range("AG7").Value = "1"
$ws.Range("AJ7").Value = "N/A"

# Copy the formatting (style index, borders, etc.) of row 6 onto row 7 so every
# cell in the new row carries the same cell style ("s=2") the other data rows
# use -- this also stamps the still-empty cells with that style without
# touching the values already written above.
$ws.Range("A6:AJ6").Copy()
$ws.Range("A7:AJ7").PasteSpecial(-4122)

# Row 7 is a detail/slice row like rows 2-6 and stays hidden behind the filter.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(5).Hidden = $true
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(7).Hidden = $true

# ---------------------------------------------------------------------------
# 3. Column A width shrinks now that "Extension.valueCodeableConcept" (30
#    characters) is gone and "Extension.extension" (19 characters) is the
#    longest remaining path in the column.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.17

# ---------------------------------------------------------------------------
# 4. AutoFilter / _FilterDatabase / conditional formatting all need to grow
#    from row 6 to row 7.
# ---------------------------------------------------------------------------
$ws.Range("A1:AJ7").AutoFilter()
$ws.Range("A1:AJ7").AutoFilter()
$ws.Range("A1:AJ7").AutoFilter(7, "<> ", 1)
$ws.Range("A1:AJ7").AutoFilter(27, @(""), 7)

$wb.Names.Item(1).RefersTo = "=Elements!`$A`$1:`$AJ`$7"

$ws.Range("A2:AI5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A2:AI6"))
